$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D4').Value = "1) Link user to respective page"
$ws.Range('B6').Value = "Browse tab link (centre of content, browse button)"
$ws.Range('D6').Value = "1) Users able to look for soap collection images, soap label name,color, country and for skin type`n2) Users able to look at entire collection at browse page `n3) Users able to filter and search for wanted soap collection`n*Refer to User at Browse Page for more information(Test Case 4)"
$ws.Range('B8').Value = "Browse or Add button link"
$ws.Range('C8').Value = "1) Clicking will link user to respective page"
$ws.Range('D8').Value = "1) Users able to browse search for collections, look at entire collection at browse page`n*Refer to User at Browse Page for more information(Test Case 4)`n2) Users able to add new collections at add page`n*Refer to User at Add Page for more information(Test Case 5)"
$ws.Range('D10').Value = "1) Clicking on dropdown will enable users to input search for following fields: color, country, max & min cost, skin type.`n"
$ws.Range('D13').Value = "`n1) Results of relavant collections will display after the accordian box upon clicking on the search button"
$ws.Range('C15').Value = "2)  Clicking on edit "
$ws.Range('C16').Value = "3) Clicking on more "
$ws.Range('C17').Value = "4) Clicking on delete cross icon at  top right most corner"
$ws.Range('D17').Value = "1) Bring users to confirm delete modal page, prompt users to be sure to delete`n2) Users able to delete by clicking on the delete button at the bottom of modal and item will be removed from collection  "

$ws.Rows.Item(6).RowHeight = 129
$ws.Rows.Item(8).RowHeight = 123.6
$ws.Rows.Item(13).RowHeight = 57.6

$ws.Columns.Item(3).ColumnWidth = 66.5546875

$ws.PageSetup.Zoom = 59

$ws.Range('D8').Select()
